# Clases.xlsx edit: "lucia cambia xls de clases"
#
# Applies (to the extent the COM surface allows):
#   1. Rewrites the "Barabas 2010 [pdf](...)" reading link in C2 to a direct
#      Google-Drive download URL and clears the trailing "Aronow ... [pdf](url)"
#      link target.
#   2. Bumps the zoom level of the sheet view from 90% to 140%.
#   3. Moves the active selection from C3 to C2.
#   4. Grows row 2's height from 64 to 80 points (to fit the edited, wrapped text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the shared-string text used by C2 ---------------------------
# Keep the first line (with its non-breaking space before "&") untouched and
# only change the second line, per the authored diff.
$nbsp = [char]0x00A0
$line1 = "Gerber " + $nbsp + "& Green 2012. FEDAI [Descarga] https://drive.google.com/drive/folders/14HDw0lx7v8cduNtj2XNvvZ5fm_lQ7Z6y?usp=sharing)"
$line2 = "Barabas 2010 [pdf](https://drive.google.com/u/0/uc?id=15SqCaheQIA_Eg8Q6CxkkF5Gdt2dPdK1Y&export=download)  Aronow et al 2015 [pdf]()"
$ws.Range("C2").Value = $line1 + [char]10 + $line2

# --- 2. Zoom the sheet view to 140% -----------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 140

# --- 3. Move the selection to C2 --------------------------------------------
$ws.Range("C2").Select()

# --- 4. Resize row 2 to fit the updated wrapped text ------------------------
$ws.Rows.Item(2).RowHeight = 80
